# Scoreboard functionality with reset
# - Adds a "PLAY AGAIN" button (duplicate of the existing "PLAY GAME" button)
# - Refreshes the cached footer date placeholder text across every layout/master
# - Normalizes the Articulate project tag ordering (SLIDE_COUNT, then PROJECT_OPEN)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Add the "PLAY AGAIN" button to slide 1, mirroring the "PLAY GAME" button.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

$playGame = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText -and $candidate.TextFrame.TextRange.Text -eq "PLAY GAME") {
        $playGame = $candidate
    }
}

if ($playGame -ne $null) {
    $newRange = $playGame.Duplicate()
    $playAgain = $newRange.Item(1)
    $playAgain.Name = "Rectangle: Rounded Corners 2"
    # Left/Top/Width/Height are expressed in points; convert the target EMU
    # geometry (below the PLAY GAME button) down to points for an exact placement.
    $playAgain.Left = 2192694 / 12700.0
    $playAgain.Top = 5460896 / 12700.0
    $playAgain.Width = 1894114 / 12700.0
    $playAgain.Height = 531845 / 12700.0
    $playAgain.TextFrame.TextRange.Text = "PLAY AGAIN"
}

# ---------------------------------------------------------------------------
# 2. Refresh the cached date footer ("datetimeFigureOut") on every layout and
#    on the slide master so it reads 21/09/2022.
# ---------------------------------------------------------------------------
$newDateText = "21/09/2022"
$master = $p.SlideMaster

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shape = $layout.Shapes.Item($j)
        if ($shape.Name -like "Date Placeholder*") {
            $shape.TextFrame.TextRange.Text = $newDateText
        }
    }
}

for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $shape = $master.Shapes.Item($j)
    if ($shape.Name -like "Date Placeholder*") {
        $shape.TextFrame.TextRange.Text = $newDateText
    }
}

# ---------------------------------------------------------------------------
# 3. Reorder the project-level Articulate tags so ARTICULATE_SLIDE_COUNT is
#    listed before ARTICULATE_PROJECT_OPEN.
# ---------------------------------------------------------------------------
$tags = $p.Tags
$tags.Delete("ARTICULATE_PROJECT_OPEN")
$tags.Add("ARTICULATE_PROJECT_OPEN", "0")
